$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.773.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.411.45'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.83%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '116.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +10.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '318.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.640'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.55%  '
$ws.Range("E8").Value = '  -0.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.634'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.12'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0935'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.66%  '
$ws.Range("E12").Value = '  +6.36%  '
$ws.Range("E13").Value = '  +2.60%  '
$ws.Range("E14").Value = '  +2.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '16.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.779.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.417.62'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.736.58'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.67%  '
$ws.Range("E20").Value = '  +4.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.23'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.57'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.39'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.34%  '
$ws.Range("E26").Value = '  -0.80%  '
$ws.Range("E27").Value = '  +5.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '173.64'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.96'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +12.35%  '
$ws.Range("E35").Value = '  +11.12%  '
$ws.Range("E36").Value = '  +2.21%  '
$ws.Range("E37").Value = '  +7.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.27'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +17.58%  '
$ws.Range("E39").Value = '  +12.03%  '
$ws.Range("E40").Value = '  +5.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.79'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +13.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.85'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '117.25'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.62'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +11.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +15.36%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '40.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +13.29%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.39%  '
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0995'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +16.45%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.93'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.71%  '
$ws.Range("B42").Value = 'Celestia'
$ws.Range("C42").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.80'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +12.17%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.242'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.86%  '
$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.57'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -9.54%  '
$ws.Range("B45").Value = 'MultiversX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '72.62'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.61%  '
$ws.Range("B46").Value = 'ordi'
$ws.Range("C46").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '86.71'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +11.63%  '
$ws.Range("B47").Value = 'FirstDigitalUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.50%  '
